$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.871.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.668.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.260"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0622"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.24%  "

$ws.Range("E11").Value = "  -1.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.904.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.669.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("E14").Value = "  -1.30%  "

$ws.Range("E15").Value = "  -1.56%  "

$ws.Range("E16").Value = "  -0.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "250.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.838.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.65%  "

$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.43%  "

$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.13%  "

$ws.Range("E23").Value = "  -2.32%  "

$ws.Range("E24").Value = "  -1.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("E28").Value = "  -0.30%  "

$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("E30").Value = "  +6.20%  "

$ws.Range("E31").Value = "  +0.30%  "

$ws.Range("E32").Value = "  -0.71%  "

$ws.Range("E33").Value = "  -3.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.431.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.44%  "

$ws.Range("E35").Value = "  -5.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.931"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.584"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.81%  "

$ws.Range("E39").Value = "  -1.14%  "

$ws.Range("E40").Value = "  -2.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.812.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.788"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("E47").Value = "  +4.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0112"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.56%  "
